$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the used range dimension by writing header K1, matching header style from J1
$ws.Range("J1").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("K1").Value = "kg_quantities_per_portion"
$ws.Range("K2").Value = "[0.01, 0.016666666666666666, 0.03333333333333333, 0.00020537166666666665, 0.13333333333333333, 0.0008214866666666666, 0.015, 0.06666666666666667, 0.049999999999999996, 0.011666666666666667, 0.008333333333333333, 0.0]"
$ws.Range("K3").Value = "[0.03, 0.025, 0.05, 0.00123223, 0.02, 0.15, 0.1, 0.025, 0.025, 0.0]"
$ws.Range("K4").Value = "[0.05000000000000001, 0.03333333333333333, 0.02, 0.049999999999999996, 0.020833333333333332, 0.0]"
$ws.Range("K5").Value = "[0.05, 0.08, 0.025, 0.15000000000000002, 0.00123223, 0.0, 0.00492892, 0.0, 0.0075, 0.1, 0.14]"
$ws.Range("K6").Value = "[0.25, 0.015, 0.15000000000000002, 0.0, 0.0225, 0.0225, 0.1]"
$ws.Range("K7").Value = "[0.1, 0.0075, 0.0075, 0.1, 0.00123223, 0.05, 0.1, 0.025, 0.025, 0.05, 0.0075, 0.0]"
$ws.Range("K8").Value = "[0.00375, 0.05, 0.0075, 0.025, 0.00246446, 0.025, 0.005625, 0.25, 0.125, 0.08]"
$ws.Range("K9").Value = "[0.25, 0.03333333333333333, 0.005, 0.0025, 0.016666666666666666, 0.03333333333333333, 0.072, 0.0016429733333333332, 0.0016429733333333332, 0.0008214866666666666, 0.0008214866666666666, 0.0]"
$ws.Range("K10").Value = "[0.0062499999999999995, 0.16666666666666666, 0.016666666666666666, 0.03333333333333333, 0.03333333333333333, 0.03333333333333333, 0.0075, 0.0025, 0.024999999999999998, 0.03333333333333333, 0.016666666666666666, 0.0625]"
$ws.Range("K11").Value = "[0.015, 0.05, 0.05, 0.2, 0.05, 0.1, 0.0225, 0.00375, 0.1]"
$ws.Range("K12").Value = "[0.005625, 0.025, 0.05, 0.025, 0.000616115, 0.2, 0.025, 0.125, 0.30000000000000004, 0.0375, 0.0, 0.0]"
$ws.Range("K13").Value = "[0.1, 0.00375, 0.001875, 0.025, 0.000616115, 0.025, 0.025, 0.0125, 0.05, 0.05, 0.0, 0.0]"
$ws.Range("K14").Value = "[0.011666666666666667, 0.011666666666666667, 0.16666666666666666, 0.016666666666666666, 0.013333333333333334, 0.0016429733333333332, 0.0, 0.09999999999999999, 0.041666666666666664, 0.0, 0.049999999999999996, 0.20000000000000004, 0.008333333333333333]"
$ws.Range("K15").Value = "[0.05, 0.0125, 0.05, 0.05, 0.05, 0.0, 0.25, 0.0, 0.05]"
$ws.Range("K16").Value = "[0.08750000000000001, 0.05, 0.0075, 0.1, 0.1, 0.00375, 0.0375, 0.00123223, 0.05, 0.00123223, 0.0625]"
$ws.Range("K17").Value = "[0.08, 0.0, 0.03333333333333333, 0.04666666666666667, 0.06666666666666667, 0.04666666666666667, 0.016666666666666666, 0.0]"
$ws.Range("K18").Value = "[0.10833333333333334, 0.006666666666666667, 0.006666666666666667, 0.0016429733333333332, 0.024999999999999998, 0.03, 0.0, 0.049999999999999996, 0.049999999999999996, 0.0]"
$ws.Range("K19").Value = "[0.11666666666666668, 0.0016429733333333332, 0.0016429733333333332, 0.0075, 0.0, 0.10000000000000002, 0.005, 0.26666666666666666, 0.0032859466666666664, 0.049999999999999996, 0.13333333333333333, 0.008333333333333333, 0.0]"
$ws.Range("K20").Value = "[0.00375, 0.0075, 0.005625, 0.2, 0.025, 0.025, 0.0075, 0.03125, 0.025, 0.07500000000000001, 0.1]"
$ws.Range("K21").Value = "[0.23333333333333336, 0.01, 0.03333333333333333, 0.13333333333333333, 0.0, 0.13333333333333333, 0.0, 0.0008214866666666666, 0.09999999999999999]"
$ws.Range("K22").Value = "[0.05, 0.1, 0.0225, 0.2, 0.1, 0.0, 0.015, 0.0]"
$ws.Range("K23").Value = "[0.025, 0.0, 0.05, 0.00375, 0.0075, 0.1, 0.0075, 0.001875, 0.07500000000000001, 0.05, 0.01125, 0.075, 0.00375, 0.0, 0.0]"
$ws.Range("K24").Value = "[0.05, 0.00375, 0.015, 0.0075, 0.07500000000000001, 0.025, 0.00375, 0.075, 0.0125, 0.00375, 0.1, 0.025]"
$ws.Range("K25").Value = "[0.075, 0.075, 0.00375, 0.0, 0.0, 0.0075, 0.00375, 0.0125, 0.05]"
